# repull data, push all data, mean calculation
# Update the dSF (column F) values for several rows to reflect repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = -1
    7  = -5
    8  = -3
    11 = -3
    15 = -5
    17 = -1
    18 = -3
    19 = -3
    20 = -1
    21 = -1
    22 = -2
    23 = 1
    24 = 1
    25 = -1
    26 = 2
    27 = -2
    29 = -2
    30 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
